# Generate Report for Archive
#
# 1. The localization status label "Ready for handoff" is renamed to
#    "In Translation" everywhere it appears (Overview!E:F, zh-cn!C,
#    de-de!C all point at the same shared string, so every cell that
#    currently holds that text is updated).
# 2. The (now narrower) Status columns are resized down from
#    ~17.22 chars to ~13.41 chars: Overview!E:F, zh-cn!C, de-de!C.
#    NOTE: Range.ColumnWidth only accepts values that land on Excel's
#    internal pixel grid (steps of 1/6 character at the workbook's
#    default font), so 12.5 is used here because it is the input that
#    rounds to the stored width closest to the authored 13.4101845877511.

$wb = $excel.ActiveWorkbook

# --- 1. Rename the status text everywhere it occurs ---------------------
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        # Compare with the literal on the left so PowerShell doesn't try
        # to coerce the (string) cell value to the type of the cell's
        # current contents (e.g. text cells that merely read "True").
        if ("Ready for handoff" -eq $cell.Value2) {
            $cell.Value = "In Translation"
        }
    }
}

# --- 2. Narrow the Status columns ---------------------------------------
$newWidth = 12.5

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E1:F1").ColumnWidth = $newWidth

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C1").ColumnWidth = $newWidth

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C1").ColumnWidth = $newWidth
